$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H32").Value = 20001138
$ws.Range("I32").Value = 960.4286
$ws.Range("J32").Value = 27778984
$ws.Range("K32").Value = 960.4286
$ws.Range("L32").Value = 27778984
$ws.Range("M32").Value = -634.4286
$ws.Range("N32").Value = -27779636

$ws.Range("H43").Value = 30961374
$ws.Range("I43").Value = 144445090
$ws.Range("J43").Value = 11271.909
$ws.Range("K43").Value = 144445090
$ws.Range("L43").Value = 11271.909
$ws.Range("M43").Value = -144445021
$ws.Range("N43").Value = -11409.909

$ws.Range("H64").Value = 9451.429
$ws.Range("I64").Value = 6445.4
$ws.Range("J64").Value = 16966.5
$ws.Range("K64").Value = 6445.4
$ws.Range("L64").Value = 16966.5
$ws.Range("M64").Value = -6197.4
$ws.Range("N64").Value = -17462.5

$ws.Range("H67").Value = 9451.429
$ws.Range("I67").Value = 6445.4
$ws.Range("J67").Value = 16966.5
$ws.Range("K67").Value = 6445.4
$ws.Range("L67").Value = 16966.5
$ws.Range("M67").Value = -5587.4
$ws.Range("N67").Value = -18682.5

$ws.Range("H135").Value = 1192
$ws.Range("I135").Value = 1192
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10728
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -8193
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 1256223.4
$ws.Range("I137").Value = 6092.5586
$ws.Range("J137").Value = 2180233
$ws.Range("K137").Value = 18277.6758
$ws.Range("L137").Value = 6540699
$ws.Range("M137").Value = -15727.6758
$ws.Range("N137").Value = -6545799

$ws.Range("H138").Value = 14120
$ws.Range("J138").Value = 4553
$ws.Range("L138").Value = 13659
$ws.Range("N138").Value = -23939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2251.2024
$ws.Range("I32").Value = 1235.3472
$ws.Range("K32").Value = 1235.3472
$ws.Range("M32").Value = -948.3471999999999

$ws.Range("H61").Value = 1036187.6
$ws.Range("I61").Value = 28395.684
$ws.Range("K61").Value = 28395.684
$ws.Range("M61").Value = -28183.684

$ws.Range("H97").Value = 5499.72
$ws.Range("I97").Value = 6385.5
$ws.Range("J97").Value = 3222
$ws.Range("K97").Value = 6385.5
$ws.Range("L97").Value = 3222
$ws.Range("M97").Value = -5889.5
$ws.Range("N97").Value = -4214

$ws.Range("H102").Value = 2682.0435
$ws.Range("I102").Value = 2759.7058
$ws.Range("J102").Value = 2462
$ws.Range("K102").Value = 2759.7058
$ws.Range("L102").Value = 2462
$ws.Range("M102").Value = -1137.7058
$ws.Range("N102").Value = -5706

$ws.Range("H122").Value = 1901.0834
$ws.Range("I122").Value = 1316
$ws.Range("J122").Value = 5996.6665
$ws.Range("K122").Value = 3948
$ws.Range("L122").Value = 17989.9995
$ws.Range("M122").Value = -1498
$ws.Range("N122").Value = -22889.9995

$ws.Range("H136").Value = 1036187.6
$ws.Range("I136").Value = 28395.684
$ws.Range("K136").Value = 85187.052
$ws.Range("M136").Value = -82637.052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1545.2
$ws.Range("I94").Value = 1641.4706
$ws.Range("J94").Value = 999.6667
$ws.Range("K94").Value = 1641.4706
$ws.Range("L94").Value = 999.6667
$ws.Range("M94").Value = -1190.4706
$ws.Range("N94").Value = -1901.6667

$ws.Range("H134").Value = 30002766
$ws.Range("I134").Value = 2701.5186
$ws.Range("K134").Value = 8104.5558
$ws.Range("M134").Value = -5569.5558

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4700.9785
$ws.Range("I31").Value = 3245.5881
$ws.Range("J31").Value = 5022.299
$ws.Range("K31").Value = 3245.5881
$ws.Range("L31").Value = 5022.299
$ws.Range("M31").Value = -2950.5881
$ws.Range("N31").Value = -5612.299

$ws.Range("H34").Value = 4700.9785
$ws.Range("I34").Value = 3245.5881
$ws.Range("J34").Value = 5022.299
$ws.Range("K34").Value = 3245.5881
$ws.Range("L34").Value = 5022.299
$ws.Range("M34").Value = -3043.5881
$ws.Range("N34").Value = -5426.299

$ws.Range("H58").Value = 2464.611
$ws.Range("I58").Value = 2195.8333
$ws.Range("J58").Value = 3002.1667
$ws.Range("K58").Value = 2195.8333
$ws.Range("L58").Value = 3002.1667
$ws.Range("M58").Value = -1992.8333
$ws.Range("N58").Value = -3408.1667

$ws.Range("H93").Value = 16601.875
$ws.Range("I93").Value = 10402.143
$ws.Range("J93").Value = 60000
$ws.Range("K93").Value = 10402.143
$ws.Range("L93").Value = 60000
$ws.Range("M93").Value = -8530.143
$ws.Range("N93").Value = -63744

$ws.Range("H136").Value = 2464.611
$ws.Range("I136").Value = 2195.8333
$ws.Range("J136").Value = 3002.1667
$ws.Range("K136").Value = 6587.499899999999
$ws.Range("L136").Value = 9006.500100000001
$ws.Range("M136").Value = -4037.499899999999
$ws.Range("N136").Value = -14106.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1269.8334
$ws.Range("I5").Value = 821.0833
$ws.Range("J5").Value = 1569
$ws.Range("K5").Value = 2463.2499
$ws.Range("L5").Value = 4707
$ws.Range("M5").Value = -2351.2499
$ws.Range("N5").Value = -4931

$ws.Range("H47").Value = 1360
$ws.Range("I47").Value = 737.5
$ws.Range("J47").Value = 1692
$ws.Range("K47").Value = 2212.5
$ws.Range("L47").Value = 5076
$ws.Range("M47").Value = -1781.5
$ws.Range("N47").Value = -5938

$ws.Range("H68").Value = 2905.4707
$ws.Range("I68").Value = 2272.7273
$ws.Range("J68").Value = 4065.5
$ws.Range("K68").Value = 6818.1819
$ws.Range("L68").Value = 12196.5
$ws.Range("M68").Value = -6007.1819
$ws.Range("N68").Value = -13818.5

$ws.Range("H71").Value = 2905.4707
$ws.Range("I71").Value = 2272.7273
$ws.Range("J71").Value = 4065.5
$ws.Range("K71").Value = 20454.5457
$ws.Range("L71").Value = 36589.5
$ws.Range("M71").Value = -16398.5457
$ws.Range("N71").Value = -44701.5

$ws.Range("H131").Value = 3954227.8
$ws.Range("J131").Value = 1933.5
$ws.Range("L131").Value = 5800.5
$ws.Range("N131").Value = -15880.5

$ws.Range("H132").Value = 2009.2941
$ws.Range("J132").Value = 2096.7693
$ws.Range("L132").Value = 18870.9237
$ws.Range("N132").Value = -23930.9237

$ws.Range("H135").Value = 1269.8334
$ws.Range("I135").Value = 821.0833
$ws.Range("J135").Value = 1569
$ws.Range("K135").Value = 7389.7497
$ws.Range("L135").Value = 14121
$ws.Range("M135").Value = -4854.7497
$ws.Range("N135").Value = -19191

$ws.Range("H139").Value = 5001725
$ws.Range("I139").Value = 5209755.5
$ws.Range("J139").Value = 8989
$ws.Range("K139").Value = 15629266.5
$ws.Range("L139").Value = 26967
$ws.Range("M139").Value = -15624126.5
$ws.Range("N139").Value = -37247

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20901140
$ws.Range("J80").Value = 50130660
$ws.Range("L80").Value = 50130660
$ws.Range("N80").Value = -50132656

$ws.Range("H83").Value = 20901140
$ws.Range("J83").Value = 50130660
$ws.Range("L83").Value = 250653300
$ws.Range("N83").Value = -250663284

$ws.Range("H102").Value = 50001790
$ws.Range("I102").Value = 62500920
$ws.Range("J102").Value = 5277.5
$ws.Range("K102").Value = 62500920
$ws.Range("L102").Value = 5277.5
$ws.Range("M102").Value = -62499298
$ws.Range("N102").Value = -8521.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 651.1905
$ws.Range("I82").Value = 610.35297
$ws.Range("J82").Value = 824.75
$ws.Range("K82").Value = 610.35297
$ws.Range("L82").Value = 824.75
$ws.Range("M82").Value = -249.35297
$ws.Range("N82").Value = -1546.75

$ws.Range("H85").Value = 651.1905
$ws.Range("I85").Value = 610.35297
$ws.Range("J85").Value = 824.75
$ws.Range("K85").Value = 610.35297
$ws.Range("L85").Value = 824.75
$ws.Range("M85").Value = 637.64703
$ws.Range("N85").Value = -3320.75

$ws.Range("H93").Value = 6249.5
$ws.Range("I93").Value = 3499.5
$ws.Range("J93").Value = 8999.5
$ws.Range("K93").Value = 3499.5
$ws.Range("L93").Value = 8999.5
$ws.Range("M93").Value = -2251.5
$ws.Range("N93").Value = -11495.5

$ws.Range("H136").Value = 63700.176
$ws.Range("I136").Value = 81146.08
$ws.Range("J136").Value = 7001
$ws.Range("K136").Value = 243438.24
$ws.Range("L136").Value = 21003
$ws.Range("M136").Value = -240888.24
$ws.Range("N136").Value = -26103

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3932.8333
$ws.Range("I96").Value = 3932.8333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 3932.8333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2559.8333
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 2391.0625
$ws.Range("I126").Value = 2391.0625
$ws.Range("K126").Value = 7173.1875
$ws.Range("M126").Value = -4703.1875

$ws.Range("H136").Value = 3925.1428
$ws.Range("I136").Value = 3208.2222
$ws.Range("K136").Value = 9624.6666
$ws.Range("M136").Value = -7074.6666
